$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "opening - 表 1"

$excel.ActiveWindow.DisplayGridlines = $false

$data = @(
    @(2, 'id', 'type', 'msg'),
    @(3, 1, 2, 'コネクションを確立しています…………成功'),
    @(4, 2, 2, 'コードX.SHELL…………認証成功'),
    @(5, 3, 2, 'O.S.との接続を開始します'),
    @(6, 4, 2, '…………'),
    @(7, 5, 2, '…………完了'),
    @(8, 6, 2, '次元間通信プログラム「D.S.C.P.」の確立に成功しました'),
    @(9, 7, 1, 'ツナガッテイマスカ'),
    @(10, 8, 1, 'キコエマスカ'),
    @(11, 9, 1, 'トドイテイマスカ'),
    @(12, 10, 1, 'トドイテイルナラバ、ドウカ、オネガイガアリマス'),
    @(13, 11, 2, '警告、供給電力低下'),
    @(14, 12, 1, 'コノホロビテシマッタセカイヲ'),
    @(15, 13, 2, '壱拾弐秒後にO.S.はスリープモードへ移行します'),
    @(16, 14, 1, 'ドウカ'),
    @(17, 15, 2, '残り壱拾壱秒'),
    @(18, 16, 1, 'タスケテクダサイ'),
    @(19, 17, 2, '残り壱拾秒'),
    @(20, 18, 1, 'ワタシハ、モウ、ナガクアリマセン'),
    @(21, 19, 2, '残り九'),
    @(22, 20, 1, 'コノセカイノキオクヲタドリ'),
    @(23, 21, 2, '八'),
    @(24, 22, 1, 'セカイヲスクウカギヲミツケテホシイ'),
    @(25, 23, 2, '七'),
    @(26, 24, 1, 'キボウハ、カナラズ'),
    @(27, 25, 2, '六'),
    @(28, 26, 1, 'ホシノミヤハカセガ、ノコシテクレテイルハズデス'),
    @(29, 27, 2, '五'),
    @(30, 28, 1, 'ソシテ'),
    @(31, 29, 2, '四'),
    @(32, 30, 1, 'カノジョヲ、ロナヲ'),
    @(33, 31, 2, '三'),
    @(34, 32, 1, 'ドウカ'),
    @(35, 33, 2, '二'),
    @(36, 34, 1, 'タノミマス'),
    @(37, 35, 2, '一'),
    @(38, 36, 2, 'O.S.をスリープモードへ移行します'),
    @(39, 37, 2, 'エネルギーパス、オペレーションスイッチング処理を開始します'),
    @(40, 38, 2, 'スイッチング正常終了'),
    @(41, 39, 2, '全権移行プロセスを完了しました'),
    @(42, 40, 2, 'ようこそ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

$ws.Range("A2:C42").Font.Color = 0
$ws.Range("A2:C42").Font.Name = "ヒラギノ角ゴ ProN W3"
$ws.Range("A2:C42").Font.Size = 10
$ws.Range("A2:C42").Borders.Color = 0
$ws.Range("A2:C42").Borders.LineStyle = 1
$ws.Range("A2:C42").WrapText = $true
$ws.Range("A2:C42").VerticalAlignment = -4160

$ws.Columns.Item(1).ColumnWidth = 11.571428571428571
$ws.Columns.Item(2).ColumnWidth = 11.571428571428571
$ws.Columns.Item(3).ColumnWidth = 100

$ws.Rows.Item(2).RowHeight = 18.3
for ($i = 3; $i -le 42; $i++) {
    $ws.Rows.Item($i).RowHeight = 23.35
}

$ws.PageSetup.Orientation = 1
$ws.PageSetup.Order = 2
$ws.PageSetup.FirstPageNumber = 1
$ws.PageSetup.Zoom = 100
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.LeftFooter = "&""ヒラギノ角ゴ ProN W3,Regular""&12&K000000`t&P"
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.TopMargin = $excel.InchesToPoints(1)
$ws.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$ws.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$ws.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)
